$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "B2" = 0.04642948062583634
    "C2" = 0.2965858768643577
    "D2" = 0.1211486527137855
    "E2" = 0.3480641502852391
    "F2" = 0.3579752231001048
    "B3" = 0.07710192660348104
    "C3" = 0.3275181423519703
    "D3" = 0.1806004724994779
    "E3" = 0.4249711431373639
    "F3" = 0.4349832324714388
    "B4" = 0.06899662371576315
    "C4" = 0.2620176638063245
    "D4" = 0.09766707149571256
    "E4" = 0.3125173139135055
    "F4" = 0.3183592094900543
    "B5" = 0.1035751530193068
    "C5" = 0.3321382172543214
    "D5" = 0.1880858394500914
    "E5" = 0.4336886434414572
    "F5" = 0.4416942719119532
    "B6" = 0.1569404967837869
    "C6" = 0.3655869259556198
    "D6" = 0.2335277915346266
    "E6" = 0.4832471329812797
    "F6" = 0.481776195163832
    "B7" = 0.1649274464288797
    "C7" = 0.4479266052260063
    "D7" = 0.276931746074847
    "E7" = 0.5262430484812574
    "F7" = 0.5300443556207896
    "B8" = 0.2184461953069695
    "C8" = 0.4442913753428973
    "D8" = 0.3428095254168975
    "E8" = 0.5854993812267418
    "F8" = 0.5950705354891705
    "B9" = 0.04181093721508944
    "C9" = 0.2983623655071265
    "D9" = 0.09486659670334104
    "E9" = 0.308004215398655
    "F9" = 0.3737347499882843
    "B10" = -0.2839541438535775
    "C10" = 0.2839541438535775
    "D10" = 0.08062995581161821
    "E10" = 0.2839541438535775
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
